$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Header text changes (row 1)
$ws.Range("L1").Value = "% Not in need"
$ws.Range("M1").Value = "# Not in need"

# Row 7 (Girls -> Girls (5-17 y.o.))
$ws.Range("A7").Value = "Girls (5-17 y.o.)"
$ws.Range("D7").Value = 26.8
$ws.Range("E7").Value = 335825
$ws.Range("F7").Value = 4.4
$ws.Range("G7").Value = 54947
$ws.Range("H7").Value = 3.9
$ws.Range("I7").Value = 48329
$ws.Range("J7").Value = 1.3
$ws.Range("K7").Value = 15870
$ws.Range("L7").Value = 63.6
$ws.Range("M7").Value = 796577
$ws.Range("N7").Value = 36.4
$ws.Range("O7").Value = 454972

# Row 8 (Boys -> Boys (5-17 y.o.))
$ws.Range("A8").Value = "Boys (5-17 y.o.)"
$ws.Range("D8").Value = 28.6
$ws.Range("E8").Value = 451938
$ws.Range("F8").Value = 3.6
$ws.Range("G8").Value = 56971
$ws.Range("H8").Value = 4.1
$ws.Range("I8").Value = 65185
$ws.Range("J8").Value = 1.6
$ws.Range("K8").Value = 24672
$ws.Range("L8").Value = 62.2
$ws.Range("M8").Value = 983644
$ws.Range("N8").Value = 37.8
$ws.Range("O8").Value = 598765

# Row 9 (ECE (5 y.o.))
$ws.Range("D9").Value = 26.8
$ws.Range("E9").Value = 77773
$ws.Range("F9").Value = 0.8
$ws.Range("G9").Value = 2364
$ws.Range("H9").Value = 2.4
$ws.Range("I9").Value = 6904
$ws.Range("J9").Value = 1.5
$ws.Range("K9").Value = 4234
$ws.Range("L9").Value = 68.59999999999999
$ws.Range("M9").Value = 198957
$ws.Range("N9").Value = 31.4
$ws.Range("O9").Value = 91274

# Row 10 (Primary school)
$ws.Range("D10").Value = 22.2
$ws.Range("E10").Value = 242310
$ws.Range("F10").Value = 0.9
$ws.Range("G10").Value = 10040
$ws.Range("H10").Value = 5.1
$ws.Range("I10").Value = 55205
$ws.Range("J10").Value = 1.8
$ws.Range("K10").Value = 19733
$ws.Range("L10").Value = 70
$ws.Range("M10").Value = 762696
$ws.Range("N10").Value = 30
$ws.Range("O10").Value = 327288

# Row 11 (Upper primary school -> Intermediate school-level)
$ws.Range("A11").Value = "Intermediate school-level"
$ws.Range("D11").Value = 26.4
$ws.Range("E11").Value = 229861
$ws.Range("F11").Value = 3
$ws.Range("G11").Value = 26406
$ws.Range("H11").Value = 3.3
$ws.Range("I11").Value = 29001
$ws.Range("J11").Value = 1.4
$ws.Range("K11").Value = 12007
$ws.Range("L11").Value = 65.90000000000001
$ws.Range("M11").Value = 574713
$ws.Range("N11").Value = 34.1
$ws.Range("O11").Value = 297274

# Row 12 (Secondary school)
$ws.Range("D12").Value = 40.5
$ws.Range("E12").Value = 264609
$ws.Range("F12").Value = 12
$ws.Range("G12").Value = 78244
$ws.Range("H12").Value = 3.1
$ws.Range("I12").Value = 20489
$ws.Range("J12").Value = 0.7
$ws.Range("K12").Value = 4536
$ws.Range("L12").Value = 43.7
$ws.Range("M12").Value = 286113
$ws.Range("N12").Value = 56.3
$ws.Range("O12").Value = 367878

# Row 13 (Children with disability) is removed entirely
$ws.Rows("13:13").Delete()
